$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1051824044"
$ws.Range("D16").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E16").Value = "2201"
$ws.Range("F16").Value = 30284
$ws.Range("C17").Value = "1051824044"
$ws.Range("D17").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E17").Value = "2112"
$ws.Range("F17").Value = 36341
$ws.Range("C18").Value = "1051824044"
$ws.Range("D18").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 36341
$ws.Range("C19").Value = "1051824044"
$ws.Range("D19").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E19").Value = "2110"
$ws.Range("F19").Value = 36341
$ws.Range("C20").Value = "1051824044"
$ws.Range("D20").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E20").Value = "2109"
$ws.Range("F20").Value = 36341
$ws.Range("C21").Value = "1051824044"
$ws.Range("D21").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E21").Value = "2108"
$ws.Range("F21").Value = 36341
$ws.Range("C22").Value = "1051824044"
$ws.Range("D22").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E22").Value = "2106"
$ws.Range("F22").Value = 36341
$ws.Range("C23").Value = "1051824044"
$ws.Range("D23").Value = "LUIS JAVIER ARRIETA YEPEZ"
$ws.Range("E23").Value = "2102"
$ws.Range("F23").Value = 15748
$ws.Range("C24").Value = "1143348214"
$ws.Range("D24").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E24").Value = "2201"
$ws.Range("F24").Value = 29260
$ws.Range("C25").Value = "1143348214"
$ws.Range("D25").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E25").Value = "2112"
$ws.Range("F25").Value = 36341
$ws.Range("C26").Value = "1143348214"
$ws.Range("D26").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E26").Value = "2111"
$ws.Range("F26").Value = 36341
$ws.Range("C27").Value = "1143348214"
$ws.Range("D27").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E27").Value = "2110"
$ws.Range("F27").Value = 36341
$ws.Range("C28").Value = "1143348214"
$ws.Range("D28").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E28").Value = "2109"
$ws.Range("F28").Value = 36341
$ws.Range("C29").Value = "1143348214"
$ws.Range("D29").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E29").Value = "2108"
$ws.Range("F29").Value = 36341
$ws.Range("C30").Value = "1143348214"
$ws.Range("D30").Value = "WILBER RAFAEL ARRIETA YEPEZ"
$ws.Range("E30").Value = "2106"
$ws.Range("F30").Value = 36341
